# Adds one new weekly price record for "Jengibre" (Terminal La Palmera de
# La Serena) as row 53, pushing the existing rows 53-74 down to 54-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; Excel shifts rows 53:74 down to 54:75
# and copies formatting (incl. the date NumberFormat on column D) from the
# row above, same as inserting a row manually in the UI.
$ws.Rows("53:53").Insert()

$newRow = 53

$ws.Cells.Item($newRow, 1).Value  = 8
$ws.Cells.Item($newRow, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item($newRow, 3).Value  = 'Coquimbo'
$ws.Cells.Item($newRow, 4).Value  = 44873
$ws.Cells.Item($newRow, 5).Value  = 4
$ws.Cells.Item($newRow, 6).Value  = 100114007
$ws.Cells.Item($newRow, 7).Value  = 'Jengibre'
$ws.Cells.Item($newRow, 8).Value  = 'Sin especificar'
$ws.Cells.Item($newRow, 9).Value  = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 480
$ws.Cells.Item($newRow, 11).Value = 14000
$ws.Cells.Item($newRow, 12).Value = 15000
$ws.Cells.Item($newRow, 13).Value = 14500
$ws.Cells.Item($newRow, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item($newRow, 15).Value = 'Perú'
$ws.Cells.Item($newRow, 16).Value = 1115
$ws.Cells.Item($newRow, 17).Value = 13
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
